$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 52.55876
$ws.Range("H2").Value = 157.67628
$ws.Range("I2").Value = 0.4767987874074868
$ws.Range("J2").Value = 0.4767987874074869
$ws.Range("M2").Value = 40.25420133333333
$ws.Range("N2").Value = 120.762604
$ws.Range("O2").Value = 0.4854671023051697
$ws.Range("P2").Value = 0.4854671023051695
$ws.Range("Q2").Value = 2115.710906870346
$ws.Range("R2").Value = 19041.39816183312
$ws.Range("S2").Value = 0.2314701257053312
$ws.Range("T2").Value = 0.2314701257053312

# Row 3
$ws.Range("G3").Value = 52.55876
$ws.Range("H3").Value = 157.67628
$ws.Range("I3").Value = 0.4767987874074868
$ws.Range("J3").Value = 0.4767987874074869
$ws.Range("O3").Value = 0.05577747182450057
$ws.Range("P3").Value = 0.05577747182450056
$ws.Range("Q3").Value = 243.08342405984
$ws.Range("R3").Value = 2187.75081653856
$ws.Range("S3").Value = 0.02659463093057713
$ws.Range("T3").Value = 0.02659463093057713

# Row 4
$ws.Range("G4").Value = 52.55876
$ws.Range("H4").Value = 157.67628
$ws.Range("I4").Value = 0.4767987874074868
$ws.Range("J4").Value = 0.4767987874074869
$ws.Range("M4").Value = 36.68940733333334
$ws.Range("N4").Value = 110.068222
$ws.Range("O4").Value = 0.4424755596543956
$ws.Range("P4").Value = 0.4424755596543954
$ws.Range("Q4").Value = 1928.349754574907
$ws.Range("R4").Value = 17355.14779117416
$ws.Range("S4").Value = 0.2109718103006649
$ws.Range("T4").Value = 0.2109718103006649

# Row 5
$ws.Range("G5").Value = 52.55876
$ws.Range("H5").Value = 157.67628
$ws.Range("I5").Value = 0.4767987874074868
$ws.Range("J5").Value = 0.4767987874074869
$ws.Range("M5").Value = 1.349902
$ws.Range("N5").Value = 4.049706
$ws.Range("O5").Value = 0.01627986621593436
$ws.Range("P5").Value = 0.01627986621593436
$ws.Range("Q5").Value = 70.94917524152001
$ws.Range("R5").Value = 638.54257717368
$ws.Range("S5").Value = 0.007762220470913614
$ws.Range("T5").Value = 0.007762220470913614

# Row 6
$ws.Range("I6").Value = 0.03798452361347728
$ws.Range("J6").Value = 0.03798452361347729
$ws.Range("M6").Value = 40.25420133333333
$ws.Range("N6").Value = 120.762604
$ws.Range("O6").Value = 0.4854671023051697
$ws.Range("P6").Value = 0.4854671023051695
$ws.Range("Q6").Value = 168.5496545372426
$ws.Range("R6").Value = 1516.946890835184
$ws.Range("S6").Value = 0.0184402366110771
$ws.Range("T6").Value = 0.01844023661107711

# Row 7
$ws.Range("I7").Value = 0.03798452361347728
$ws.Range("J7").Value = 0.03798452361347729
$ws.Range("O7").Value = 0.05577747182450057
$ws.Range("P7").Value = 0.05577747182450056
$ws.Range("S7").Value = 0.002118680695617805
$ws.Range("T7").Value = 0.002118680695617806

# Row 8
$ws.Range("I8").Value = 0.03798452361347728
$ws.Range("J8").Value = 0.03798452361347729
$ws.Range("M8").Value = 36.68940733333334
$ws.Range("N8").Value = 110.068222
$ws.Range("O8").Value = 0.4424755596543956
$ws.Range("P8").Value = 0.4424755596543954
$ws.Range("Q8").Value = 153.6233915064347
$ws.Range("R8").Value = 1382.610523557912
$ws.Range("S8").Value = 0.01680722334407896
$ws.Range("T8").Value = 0.01680722334407896

# Row 9
$ws.Range("I9").Value = 0.03798452361347728
$ws.Range("J9").Value = 0.03798452361347729
$ws.Range("M9").Value = 1.349902
$ws.Range("N9").Value = 4.049706
$ws.Range("O9").Value = 0.01627986621593436
$ws.Range("P9").Value = 0.01627986621593436
$ws.Range("Q9").Value = 5.652217861064001
$ws.Range("R9").Value = 50.86996074957601
$ws.Range("S9").Value = 0.0006183829627034097
$ws.Range("T9").Value = 0.0006183829627034097

# Row 10
$ws.Range("G10").Value = 53.437349
$ws.Range("H10").Value = 160.312047
$ws.Range("I10").Value = 0.484769108051078
$ws.Range("J10").Value = 0.4847691080510781
$ws.Range("M10").Value = 40.25420133333333
$ws.Range("N10").Value = 120.762604
$ws.Range("O10").Value = 0.4854671023051697
$ws.Range("P10").Value = 0.4854671023051695
$ws.Range("Q10").Value = 2151.077805365599
$ws.Range("R10").Value = 19359.70024829038
$ws.Range("S10").Value = 0.2353394541726185
$ws.Range("T10").Value = 0.2353394541726185

# Row 11
$ws.Range("G11").Value = 53.437349
$ws.Range("H11").Value = 160.312047
$ws.Range("I11").Value = 0.484769108051078
$ws.Range("J11").Value = 0.4847691080510781
$ws.Range("O11").Value = 0.05577747182450057
$ws.Range("P11").Value = 0.05577747182450056
$ws.Range("Q11").Value = 247.1468841274161
$ws.Range("R11").Value = 2224.321957146744
$ws.Range("S11").Value = 0.02703919526570728
$ws.Range("T11").Value = 0.02703919526570728

# Row 12
$ws.Range("G12").Value = 53.437349
$ws.Range("H12").Value = 160.312047
$ws.Range("I12").Value = 0.484769108051078
$ws.Range("J12").Value = 0.4847691080510781
$ws.Range("M12").Value = 36.68940733333334
$ws.Range("N12").Value = 110.068222
$ws.Range("O12").Value = 0.4424755596543956
$ws.Range("P12").Value = 0.4424755596543954
$ws.Range("Q12").Value = 1960.584664274493
$ws.Range("R12").Value = 17645.26197847044
$ws.Range("S12").Value = 0.2144984823880629
$ws.Range("T12").Value = 0.2144984823880629

# Row 13
$ws.Range("G13").Value = 53.437349
$ws.Range("H13").Value = 160.312047
$ws.Range("I13").Value = 0.484769108051078
$ws.Range("J13").Value = 0.4847691080510781
$ws.Range("M13").Value = 1.349902
$ws.Range("N13").Value = 4.049706
$ws.Range("O13").Value = 0.01627986621593436
$ws.Range("P13").Value = 0.01627986621593436
$ws.Range("Q13").Value = 72.13518428979802
$ws.Range("R13").Value = 649.2166586081821
$ws.Range("S13").Value = 0.007891976224689378
$ws.Range("T13").Value = 0.00789197622468938

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.049338
$ws.Range("H14").Value = 0.148014
$ws.Range("I14").Value = 0.0004475809279577863
$ws.Range("J14").Value = 0.0004475809279577865
$ws.Range("M14").Value = 40.25420133333333
$ws.Range("N14").Value = 120.762604
$ws.Range("O14").Value = 0.4854671023051697
$ws.Range("P14").Value = 0.4854671023051695
$ws.Range("Q14").Value = 1.986061785384
$ws.Range("R14").Value = 17.874556068456
$ws.Range("S14").Value = 0.0002172858161427254
$ws.Range("T14").Value = 0.0002172858161427255

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.049338
$ws.Range("H15").Value = 0.148014
$ws.Range("I15").Value = 0.0004475809279577863
$ws.Range("J15").Value = 0.0004475809279577865
$ws.Range("O15").Value = 0.05577747182450057
$ws.Range("P15").Value = 0.05577747182450056
$ws.Range("Q15").Value = 0.228187460592
$ws.Range("R15").Value = 2.053687145328
$ws.Range("S15").Value = 0.00002496493259834925
$ws.Range("T15").Value = 0.00002496493259834925

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.049338
$ws.Range("H16").Value = 0.148014
$ws.Range("I16").Value = 0.0004475809279577863
$ws.Range("J16").Value = 0.0004475809279577865
$ws.Range("M16").Value = 36.68940733333334
$ws.Range("N16").Value = 110.068222
$ws.Range("O16").Value = 0.4424755596543956
$ws.Range("P16").Value = 0.4424755596543954
$ws.Range("Q16").Value = 1.810181979012
$ws.Range("R16").Value = 16.291637811108
$ws.Range("S16").Value = 0.0001980436215887552
$ws.Range("T16").Value = 0.0001980436215887552

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.049338
$ws.Range("H17").Value = 0.148014
$ws.Range("I17").Value = 0.0004475809279577863
$ws.Range("J17").Value = 0.0004475809279577865
$ws.Range("M17").Value = 1.349902
$ws.Range("N17").Value = 4.049706
$ws.Range("O17").Value = 0.01627986621593436
$ws.Range("P17").Value = 0.01627986621593436
$ws.Range("Q17").Value = 0.06660146487600001
$ws.Range("R17").Value = 0.5994131838840001
$ws.Range("S17").Value = 0.000007286557627956517
$ws.Range("T17").Value = 0.000007286557627956518
